$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.808.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.624.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.518"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.80%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.19"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0607"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.91%  "
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.855.73"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.607.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.35%  "
$ws.Range("E14").Value = "  -1.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.556"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.89"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.827.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0716"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.02%  "
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("E22").Value = "  -0.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.995"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  -0.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0480"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.08"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.402.36"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  +2.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.997"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("E37").Value = "  -1.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0169"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.82%  "
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.844"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.995"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("E45").Value = "  -1.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.764.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("E47").Value = "  -3.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.21%  "
$ws.Range("E49").Value = "  +0.76%  "
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("E51").Value = "  +0.69%  "
